$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Apply red font colour (FF0000) to the two bullet paragraphs that
#    were still missing it:
#    - "Para cada nota es posible borrarla pulsando en el boton de la derecha."
#    - "Al crearse cada nota debe realizarse una animacion."
# ------------------------------------------------------------------
$pBorrar = $d.Paragraphs(12)
$pBorrar.Range.Font.Color = 255

$pAnimacion = $d.Paragraphs(14)
$pAnimacion.Range.Font.Color = 255

# ------------------------------------------------------------------
# 2. Move the "_GoBack" bookmark from its current position (in the
#    middle of the "menu superior" paragraph) to the end of the
#    "Debe aplicarse el estilo..." paragraph, as the last edited spot.
# ------------------------------------------------------------------
try {
    $old = $d.Bookmarks("_GoBack")
    $old.Delete()
} catch {
}

$pEstilo = $d.Paragraphs(13)
$rEstilo = $pEstilo.Range
$rEstilo.MoveEnd(1, -1)
$d.Bookmarks.Add("_GoBack", $rEstilo)

# ------------------------------------------------------------------
# 3. Re-balance the "vue-cli" split: the word "vue" stays wrapped by
#    proofErr spell-check markers while "-cli" moves to the following
#    run.
# ------------------------------------------------------------------
$d.Content.Find.Execute("vue-cli", $true, $false, $false, $false, $false, $true, 1, $false, "vue", 2)
$d.Content.Find.Execute(" y adecuando la organización del código y componentes a esta plantilla.", $true, $false, $false, $false, $false, $true, 1, $false, "-cli y adecuando la organización del código y componentes a esta plantilla.", 2)

# ------------------------------------------------------------------
# 4. Re-join the "compon" / "ente" split left behind by the old
#    bookmark location into a single run.
# ------------------------------------------------------------------
$d.Content.Find.Execute("componente que maquete de un modo sencillo el tiempo actual, obtenido de alguna API REST para la ciudad desde la que se carga la web. Esta información se refresca cada 10 minutos.", $true, $false, $false, $false, $false, $true, 1, $false, "componente que maquete de un modo sencillo el tiempo actual, obtenido de alguna API REST para la ciudad desde la que se carga la web. Esta información se refresca cada 10 minutos.", 2)
